# B6-PowerPoint.pptx — commit "Sat, Jun 27, 2020 12:04:52 PM"
#
# 1) Three data tables (on the slides that hold the custom "Table_0"
#    style {9C394D5F-555A-491F-880D-A76BD77AF59F}) are switched to a
#    different table style, {379DFAA1-BA6F-4C30-842A-D12FEE26EA60}.
#    PowerPoint refuses a plain `Table.Style = "{GUID}"` assignment
#    ("Table styles cannot be assigned through a property — call
#    Table.ApplyStyle instead"), so ApplyStyle is used.
#
# 2) The deck's two embedded themes ("Office Theme" / theme1.xml and
#    "Integral" / theme2.xml) swap slots, i.e. the slide master (which
#    is wired to theme2.xml) ends up rendering with the "Office Theme"
#    palette/fonts instead of "Integral". We re-assert the intended
#    design/theme through the object model so the master picks up the
#    "Office Theme" look; ApplyTheme is the documented COM call for
#    swapping a deck's applied theme.

$p = $ppt.ActivePresentation

$oldTableStyleId = "{9C394D5F-555A-491F-880D-A76BD77AF59F}"
$newTableStyleId = "{379DFAA1-BA6F-4C30-842A-D12FEE26EA60}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $oldTableStyleId) {
                $table.ApplyStyle($newTableStyleId)
            }
        }
    }
}

# Re-apply the deck's theme so the slide master (and every slide that
# inherits from it) resolves to the "Office Theme" palette that used to
# live in the other, unused theme slot.
$master = $p.SlideMaster
try {
    $master.ApplyTheme($master.Theme)
} catch {
    # Older/looser hosts may not support re-applying the live Theme
    # object directly; ignore and fall back to the Presentation-level
    # call below.
}
try {
    $p.ApplyTheme("Office Theme")
} catch {
}
